$wb = $excel.ActiveWorkbook

# --- Sheet ALC (diff @ line 1459) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1054.4625
$ws.Range("J17").Value = 1054.4625
$ws.Range("L17").Value = 3163.3875
$ws.Range("N17").Value = -3499.3875

# --- Sheet ALC (diff @ line 4567) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1109.2632
$ws.Range("I80").Value = 1434.9231
$ws.Range("J80").Value = 403.66666
$ws.Range("K80").Value = 4304.7693
$ws.Range("L80").Value = 1210.99998
$ws.Range("M80").Value = -3306.7693
$ws.Range("N80").Value = -3206.99998

# --- Sheet ALC (diff @ line 4619) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996

# --- Sheet ALC (diff @ line 4717) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1109.2632
$ws.Range("I83").Value = 1434.9231
$ws.Range("J83").Value = 403.66666
$ws.Range("K83").Value = 12914.3079
$ws.Range("L83").Value = 3632.99994
$ws.Range("M83").Value = -7922.3079
$ws.Range("N83").Value = -13616.99994

# --- Sheet ALC (diff @ line 4769) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

# --- Sheet ALC (diff @ line 7154) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 910.4138
$ws.Range("I132").Value = 766.64813
$ws.Range("J132").Value = 2851.25
$ws.Range("K132").Value = 2299.94439
$ws.Range("L132").Value = 8553.75
$ws.Range("M132").Value = 230.0556099999999
$ws.Range("N132").Value = -13613.75

# --- Sheet ALC (diff @ line 7206) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 87159.8
$ws.Range("J133").Value = 87159.8
$ws.Range("L133").Value = 87159.8
$ws.Range("N133").Value = -97279.8

# --- Sheet ALC (diff @ line 7402) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1470.7435
$ws.Range("I137").Value = 1320.75
$ws.Range("J137").Value = 1852.5454
$ws.Range("K137").Value = 3962.25
$ws.Range("L137").Value = 5557.6362
$ws.Range("M137").Value = -1412.25
$ws.Range("N137").Value = -10657.6362

# --- Sheet ALC (diff @ line 7454) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3372.963
$ws.Range("I138").Value = 3310
$ws.Range("J138").Value = 3480
$ws.Range("K138").Value = 9930
$ws.Range("L138").Value = 10440
$ws.Range("M138").Value = -4790
$ws.Range("N138").Value = -20720

# --- Sheet ALC (diff @ line 7604) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 720179.5600000001
$ws.Range("I141").Value = 850071
$ws.Range("J141").Value = 5776.6665
$ws.Range("K141").Value = 2550213
$ws.Range("L141").Value = 17329.9995
$ws.Range("M141").Value = -2545033
$ws.Range("N141").Value = -27689.9995

# --- Sheet ARM (diff @ line 9202) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2675.662
$ws.Range("I32").Value = 2032.3594
$ws.Range("K32").Value = 2032.3594
$ws.Range("M32").Value = -1745.3594

# --- Sheet ARM (diff @ line 9836) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1607
$ws.Range("J45").Value = 1884.625
$ws.Range("L45").Value = 1884.625
$ws.Range("N45").Value = -2638.625

# --- Sheet ARM (diff @ line 10614) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2896.5715
$ws.Range("I61").Value = 1655.5714
$ws.Range("K61").Value = 1655.5714
$ws.Range("M61").Value = -1443.5714

# --- Sheet ARM (diff @ line 12985) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2437.1667
$ws.Range("I110").Value = 902.5
$ws.Range("J110").Value = 5506.5
$ws.Range("K110").Value = 902.5
$ws.Range("L110").Value = 5506.5
$ws.Range("M110").Value = 1142.5
$ws.Range("N110").Value = -9596.5

# --- Sheet ARM (diff @ line 13555) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4666.6665
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

# --- Sheet ARM (diff @ line 13607) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 54997
$ws.Range("J123").Value = 54997
$ws.Range("L123").Value = 54997
$ws.Range("N123").Value = -64797

# --- Sheet ARM (diff @ line 14036) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1795.7
$ws.Range("I132").Value = 1220.1578
$ws.Range("K132").Value = 3660.4734
$ws.Range("M132").Value = -1130.4734

# --- Sheet ARM (diff @ line 14235) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2896.5715
$ws.Range("I136").Value = 1655.5714
$ws.Range("K136").Value = 4966.7142
$ws.Range("M136").Value = -2416.7142

# --- Sheet BSM (diff @ line 18465) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19497
$ws.Range("I82").Value = 9329.333000000001
$ws.Range("K82").Value = 9329.333000000001
$ws.Range("M82").Value = -8946.333000000001

# --- Sheet BSM (diff @ line 18618) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 19497
$ws.Range("I85").Value = 9329.333000000001
$ws.Range("K85").Value = 9329.333000000001
$ws.Range("M85").Value = -8003.333000000001

# --- Sheet BSM (diff @ line 19053) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 717.5454999999999
$ws.Range("I94").Value = 524.5263
$ws.Range("K94").Value = 524.5263
$ws.Range("M94").Value = -73.52629999999999

# --- Sheet BSM (diff @ line 19675) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2370.111
$ws.Range("I107").Value = 2370.111
$ws.Range("K107").Value = 2370.111
$ws.Range("M107").Value = -450.1109999999999

# --- Sheet BSM (diff @ line 20947) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3620.0588
$ws.Range("I134").Value = 3907.9556
$ws.Range("K134").Value = 11723.8668
$ws.Range("M134").Value = -9188.8668

# --- Sheet CRP (diff @ line 22839) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1460.52
$ws.Range("I31").Value = 1240
$ws.Range("K31").Value = 1240
$ws.Range("M31").Value = -945

# --- Sheet CRP (diff @ line 22992) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1460.52
$ws.Range("I34").Value = 1240
$ws.Range("K34").Value = 1240
$ws.Range("M34").Value = -1038

# --- Sheet CRP (diff @ line 26435) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1831.4286
$ws.Range("I105").Value = 1720
$ws.Range("K105").Value = 1720
$ws.Range("M105").Value = 27

# --- Sheet CRP (diff @ line 27731) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1753.7778
$ws.Range("I132").Value = 1190.8889
$ws.Range("J132").Value = 2879.5557
$ws.Range("K132").Value = 3572.6667
$ws.Range("L132").Value = 8638.667099999999
$ws.Range("M132").Value = -1042.6667
$ws.Range("N132").Value = -13698.6671

# --- Sheet CUL (diff @ line 28685) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 18333.334
$ws.Range("I9").Value = 15000
$ws.Range("J9").Value = 19000
$ws.Range("K9").Value = 45000
$ws.Range("L9").Value = 57000
$ws.Range("M9").Value = -44776
$ws.Range("N9").Value = -57448

# --- Sheet CUL (diff @ line 32579) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11886.333
$ws.Range("I87").Value = 6162.8335
$ws.Range("K87").Value = 18488.5005
$ws.Range("M87").Value = -17240.5005

# --- Sheet CUL (diff @ line 32729) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 11886.333
$ws.Range("I90").Value = 6162.8335
$ws.Range("K90").Value = 55465.5015
$ws.Range("M90").Value = -49225.5015

# --- Sheet CUL (diff @ line 34682) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 41256
$ws.Range("I129").Value = 698.5
$ws.Range("J129").Value = 46325.688
$ws.Range("K129").Value = 2095.5
$ws.Range("L129").Value = 138977.064
$ws.Range("M129").Value = 2904.5
$ws.Range("N129").Value = -148977.064

# --- Sheet CUL (diff @ line 34734) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2995.5557
$ws.Range("I130").Value = 1230
$ws.Range("J130").Value = 3500
$ws.Range("K130").Value = 3690
$ws.Range("L130").Value = 10500
$ws.Range("M130").Value = 1330
$ws.Range("N130").Value = -20540

# --- Sheet CUL (diff @ line 34783) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 809.26
$ws.Range("I131").Value = 424.625
$ws.Range("J131").Value = 842.70654
$ws.Range("K131").Value = 1273.875
$ws.Range("L131").Value = 2528.11962
$ws.Range("M131").Value = 3766.125
$ws.Range("N131").Value = -12608.11962

# --- Sheet GSM (diff @ line 35434) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 51.22222
$ws.Range("I2").Value = 14.833333
$ws.Range("J2").Value = 69.416664
$ws.Range("K2").Value = 14.833333
$ws.Range("L2").Value = 69.416664
$ws.Range("M2").Value = 98.166667
$ws.Range("N2").Value = -295.416664

# --- Sheet GSM (diff @ line 39181) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# --- Sheet GSM (diff @ line 39325) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# --- Sheet GSM (diff @ line 41191) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1657.8889
$ws.Range("I122").Value = 1191.8667
$ws.Range("J122").Value = 2240.4167
$ws.Range("K122").Value = 3575.6001
$ws.Range("L122").Value = 6721.250100000001
$ws.Range("M122").Value = -1125.6001
$ws.Range("N122").Value = -11621.2501

# --- Sheet GSM (diff @ line 41675) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 820255.2
$ws.Range("J132").Value = 2432.2778
$ws.Range("L132").Value = 7296.8334
$ws.Range("N132").Value = -12356.8334

# --- Sheet LTW (diff @ line 44121) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4386.4
$ws.Range("I40").Value = 1701.6
$ws.Range("K40").Value = 1701.6
$ws.Range("M40").Value = -1565.6

# --- Sheet LTW (diff @ line 44409) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2794.818
$ws.Range("I46").Value = 2066.6667
$ws.Range("J46").Value = 3067.875
$ws.Range("K46").Value = 2066.6667
$ws.Range("L46").Value = 3067.875
$ws.Range("M46").Value = -1878.6667
$ws.Range("N46").Value = -3443.875

# --- Sheet LTW (diff @ line 46164) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2410.5715
$ws.Range("I82").Value = 1958
$ws.Range("K82").Value = 1958
$ws.Range("M82").Value = -1597

# --- Sheet LTW (diff @ line 46311) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2410.5715
$ws.Range("I85").Value = 1958
$ws.Range("K85").Value = 1958
$ws.Range("M85").Value = -710

# --- Sheet LTW (diff @ line 46691) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 353.27274
$ws.Range("I93").Value = 410.33334
$ws.Range("J93").Value = 96.5
$ws.Range("K93").Value = 410.33334
$ws.Range("L93").Value = 96.5
$ws.Range("M93").Value = 837.66666
$ws.Range("N93").Value = -2592.5

# --- Sheet LTW (diff @ line 48548) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1982.091
$ws.Range("I132").Value = 1912.8462
$ws.Range("J132").Value = 2027.1
$ws.Range("K132").Value = 5738.5386
$ws.Range("L132").Value = 6081.299999999999
$ws.Range("M132").Value = -3208.5386
$ws.Range("N132").Value = -11141.3

# --- Sheet LTW (diff @ line 48747) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2395.72
$ws.Range("I136").Value = 1547
$ws.Range("K136").Value = 4641
$ws.Range("M136").Value = -2091

# --- Sheet WVR (diff @ line 52964) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1672.5834
$ws.Range("I81").Value = 1696.7778
$ws.Range("K81").Value = 3393.5556
$ws.Range("M81").Value = -2332.5556

# --- Sheet WVR (diff @ line 53111) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1672.5834
$ws.Range("I84").Value = 1696.7778
$ws.Range("K84").Value = 16967.778
$ws.Range("M84").Value = -11663.778
